# Apply updated cryptocurrency price/volume data (and two row swaps)
# Commit: "Updated cryptos list on Tue Aug 27 16:52:30 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, [string]$Addr, [string]$Text)
    $cell = $Sheet.Range($Addr)
    # Force text storage so numeric-looking strings (e.g. "549.42")
    # are not silently coerced to numbers / lose formatting, then
    # drop the temporary text number-format so the cell style index
    # stays identical to its original (unstyled) state.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# Row 2
Set-CellText $ws "D2" "61.657.78"
Set-CellText $ws "E2" "  -3.15%  "

# Row 3
Set-CellText $ws "D3" "2.568.80"
Set-CellText $ws "E3" "  -5.68%  "

# Row 4
Set-CellText $ws "E4" "  -0.04%  "

# Row 5
Set-CellText $ws "D5" "549.42"
Set-CellText $ws "E5" "  -2.12%  "

# Row 6
Set-CellText $ws "D6" "154.44"
Set-CellText $ws "E6" "  -2.30%  "

# Row 7
Set-CellText $ws "D7" "0.999"
Set-CellText $ws "E7" "  +0.07%  "

# Row 8
Set-CellText $ws "D8" "0.598"
Set-CellText $ws "E8" "  +1.02%  "

# Row 9
Set-CellText $ws "E9" "  -3.10%  "

# Row 10
Set-CellText $ws "E10" "  -1.92%  "

# Row 11
Set-CellText $ws "D11" "5.43"
Set-CellText $ws "E11" "  -2.72%  "

# Row 12
Set-CellText $ws "E12" "  -2.69%  "

# Row 13
Set-CellText $ws "D13" "3.022.01"
Set-CellText $ws "E13" "  -5.70%  "

# Row 14
Set-CellText $ws "D14" "25.44"
Set-CellText $ws "E14" "  -4.59%  "

# Row 15
Set-CellText $ws "D15" "61.556.23"
Set-CellText $ws "E15" "  -3.11%  "

# Row 16
Set-CellText $ws "D16" "0.0000144"
Set-CellText $ws "E16" "  -2.37%  "

# Row 17
Set-CellText $ws "D17" "2.574.72"
Set-CellText $ws "E17" "  -5.52%  "

# Row 18
Set-CellText $ws "D18" "11.55"
Set-CellText $ws "E18" "  -5.24%  "

# Row 19
Set-CellText $ws "D19" "4.54"
Set-CellText $ws "E19" "  -2.63%  "

# Row 20
Set-CellText $ws "D20" "337.68"
Set-CellText $ws "E20" "  -3.34%  "

# Row 21
Set-CellText $ws "D21" "6.06"

# Row 22
Set-CellText $ws "D22" "0.999"
Set-CellText $ws "E22" "  -0.14%  "

# Row 23
Set-CellText $ws "D23" "0.494"
Set-CellText $ws "E23" "  -4.05%  "

# Row 24
Set-CellText $ws "D24" "63.28"
Set-CellText $ws "E24" "  -1.35%  "

# Row 25
Set-CellText $ws "E25" "  -1.22%  "

# Row 26
Set-CellText $ws "D26" "0.998"
Set-CellText $ws "E26" "  -0.17%  "

# Row 27
Set-CellText $ws "D27" "8.11"
Set-CellText $ws "E27" "  -1.37%  "

# Row 28
Set-CellText $ws "D28" "7.40"
Set-CellText $ws "E28" "  +3.00%  "

# Row 29
Set-CellText $ws "D29" "0.0₃0837"
Set-CellText $ws "E29" "  -5.52%  "

# Row 30
Set-CellText $ws "E30" "  -1.87%  "

# Row 31
Set-CellText $ws "E31" "  -5.65%  "

# Row 32
Set-CellText $ws "E32" "  +0.10%  "

# Row 33
Set-CellText $ws "D33" "159.67"
Set-CellText $ws "E33" "  -3.16%  "

# Row 34
Set-CellText $ws "E34" "  -2.90%  "

# Row 36
Set-CellText $ws "E36" "  -4.52%  "

# Row 37
Set-CellText $ws "E37" "  +0.80%  "

# Row 38: 'Bittensor' -> 'SuiNetwork'
Set-CellText $ws "B38" "SuiNetwork"
Set-CellText $ws "C38" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText $ws "D38" "0.941"
Set-CellText $ws "E38" "  -2.08%  "

# Row 39: 'SuiNetwork' -> 'Bittensor'
Set-CellText $ws "B39" "Bittensor"
Set-CellText $ws "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText $ws "D39" "333.85"
Set-CellText $ws "E39" "  -4.31%  "

# Row 40
Set-CellText $ws "D40" "5.97"
Set-CellText $ws "E40" "  -3.11%  "

# Row 41
Set-CellText $ws "D41" "3.96"
Set-CellText $ws "E41" "  -1.32%  "

# Row 42
Set-CellText $ws "D42" "37.53"
Set-CellText $ws "E42" "  -2.06%  "

# Row 43
Set-CellText $ws "D43" "20.67"
Set-CellText $ws "E43" "  -3.59%  "

# Row 44
Set-CellText $ws "D44" "0.998"
Set-CellText $ws "E44" "  +0.05%  "

# Row 45
Set-CellText $ws "D45" "2.134.03"
Set-CellText $ws "E45" "  +0.87%  "

# Row 46
Set-CellText $ws "E46" "  -3.81%  "

# Row 47
Set-CellText $ws "E47" "  -1.36%  "

# Row 48: 'EnergySwap' -> 'Hedera'
Set-CellText $ws "B48" "Hedera"
Set-CellText $ws "C48" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText $ws "D48" "0.0547"
Set-CellText $ws "E48" "  -4.86%  "

# Row 49: 'Hedera' -> 'EnergySwap'
Set-CellText $ws "B49" "EnergySwap"
Set-CellText $ws "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D49" "19.52"
Set-CellText $ws "E49" "  -5.78%  "

# Row 50
Set-CellText $ws "D50" "0.0965"
Set-CellText $ws "E50" "  -1.87%  "

# Row 51
Set-CellText $ws "E51" "  -2.94%  "
